$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# New data rows (151-164) appended to the end of the table
$data = @(
    @(44072, 15, 2080, 85, 1873, 1330, 2),
    @(44073, 10, 2090, 85, 1884, 1333, 3),
    @(44074, 9,  2099, 85, 1898, 1337, 3),
    @(44075, 18, 2117, 85, 1918, 1344, 3),
    @(44076, 12, 2129, 86, 1936, 1356, 3),
    @(44077, 11, 2140, 86, 1945, 1365, 4),
    @(44078, 20, 2160, 86, 1952, 1387, 4),
    @(44079, 16, 2176, 86, 1959, 1397, 3),
    @(44080, 8,  2184, 86, 1970, 1408, 6),
    @(44081, 6,  2190, 87, 1983, 1445, 4),
    @(44082, 10, 2200, 87, 1991, 1443, 2),
    @(44083, 11, 2211, 87, 2002, 1447, 3),
    @(44084, 8,  2219, 87, 2015, 1454, 3),
    @(44085, 20, 2239, 87, 2021, 1466, 2)
)

$startRow = 151
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# Update the view's top-left cell and selection to reflect scrolled position
$ws.Application.GoTo($ws.Range("B163"), $true)
$ws.Range("A139").Select()
$ws.Range("B163").Select()
